# Weekly update: insert two new price records (Femacal de La Calera - Tuna)
# right before the existing row 73, shifting all subsequent rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 73:74 (everything from old row 73 onward moves down by 2)
$ws.Rows("73:74").Insert()

# New row 73 - "Primera" quality, 18kg box
$ws.Range("A73").Value = 3
$ws.Range("B73").Value = "Femacal de La Calera"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44977
$ws.Range("E73").Value = 5
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100107
$ws.Range("H73").Value = "Otros"
$ws.Range("I73").Value = 100107011
$ws.Range("J73").Value = "Tuna"
$ws.Range("K73").Value = "Sin especificar"
$ws.Range("L73").Value = "Primera"
$ws.Range("M73").Value = 60
$ws.Range("N73").Value = 14000
$ws.Range("O73").Value = 14000
$ws.Range("P73").Value = 14000
$ws.Range("Q73").Value = '$/caja 18 kilos'
$ws.Range("R73").Value = "Cabildo"
$ws.Range("S73").Value = 778
$ws.Range("T73").Value = 18

# New row 74 - "Segunda" quality, 18kg box
$ws.Range("A74").Value = 3
$ws.Range("B74").Value = "Femacal de La Calera"
$ws.Range("C74").Value = "Coquimbo"
$ws.Range("D74").Value = 44977
$ws.Range("E74").Value = 5
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100107
$ws.Range("H74").Value = "Otros"
$ws.Range("I74").Value = 100107011
$ws.Range("J74").Value = "Tuna"
$ws.Range("K74").Value = "Sin especificar"
$ws.Range("L74").Value = "Segunda"
$ws.Range("M74").Value = 50
$ws.Range("N74").Value = 12000
$ws.Range("O74").Value = 12000
$ws.Range("P74").Value = 12000
$ws.Range("Q74").Value = '$/caja 18 kilos'
$ws.Range("R74").Value = "Cabildo"
$ws.Range("S74").Value = 667
$ws.Range("T74").Value = 18
